$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37294
$ws.Range("D2").Value = 53954206
$ws.Range("C3").Value = 90152
$ws.Range("D3").Value = 132180605
$ws.Range("C4").Value = 30886
$ws.Range("D4").Value = 45744271
$ws.Range("C5").Value = 8618
$ws.Range("D5").Value = 12810235
$ws.Range("C6").Value = 1964
$ws.Range("D6").Value = 2918506
$ws.Range("C7").Value = 150
$ws.Range("D7").Value = 220093
$ws.Range("C11").Value = 40788
$ws.Range("D11").Value = 55378855
$ws.Range("C12").Value = 9542
$ws.Range("D12").Value = 13803600
$ws.Range("C13").Value = 25736
$ws.Range("D13").Value = 37746042
$ws.Range("C14").Value = 8266
$ws.Range("D14").Value = 12268263
$ws.Range("C15").Value = 2131
$ws.Range("D15").Value = 3168883
$ws.Range("C16").Value = 411
$ws.Range("D16").Value = 605623
$ws.Range("C17").Value = 31
$ws.Range("D17").Value = 46500
$ws.Range("C19").Value = 10126
$ws.Range("D19").Value = 13424014
$ws.Range("C20").Value = 13266
$ws.Range("D20").Value = 19161080
$ws.Range("C21").Value = 31423
$ws.Range("D21").Value = 46123021
$ws.Range("C22").Value = 10172
$ws.Range("D22").Value = 15122897
$ws.Range("C23").Value = 2608
$ws.Range("D23").Value = 3879563
$ws.Range("C24").Value = 499
$ws.Range("D24").Value = 742592
$ws.Range("C26").Value = 11564
$ws.Range("D26").Value = 15461876
$ws.Range("C27").Value = 7559
$ws.Range("D27").Value = 10952128
$ws.Range("C28").Value = 22304
$ws.Range("D28").Value = 32737279
$ws.Range("C29").Value = 7749
$ws.Range("D29").Value = 11531802
$ws.Range("C30").Value = 1944
$ws.Range("D30").Value = 2900500
$ws.Range("C31").Value = 360
$ws.Range("D31").Value = 537415
$ws.Range("C33").Value = 8224
$ws.Range("D33").Value = 10870393
$ws.Range("C34").Value = 3177
$ws.Range("D34").Value = 4584690
$ws.Range("C35").Value = 7701
$ws.Range("D35").Value = 11247373
$ws.Range("C36").Value = 3139
$ws.Range("D36").Value = 4651754
$ws.Range("C37").Value = 816
$ws.Range("D37").Value = 1216763
$ws.Range("C38").Value = 155
$ws.Range("D38").Value = 230732
$ws.Range("C40").Value = 2406
$ws.Range("D40").Value = 3251511
$ws.Range("C41").Value = 17048
$ws.Range("D41").Value = 24658267
$ws.Range("C42").Value = 50641
$ws.Range("D42").Value = 74253776
$ws.Range("C43").Value = 18870
$ws.Range("D43").Value = 28030734
$ws.Range("C44").Value = 5564
$ws.Range("D44").Value = 8285978
$ws.Range("C45").Value = 1179
$ws.Range("D45").Value = 1759045
$ws.Range("C49").Value = 16517
$ws.Range("D49").Value = 22011851
$ws.Range("C50").Value = 1956
$ws.Range("D50").Value = 2837316
$ws.Range("C51").Value = 6706
$ws.Range("D51").Value = 9862406
$ws.Range("C52").Value = 2305
$ws.Range("D52").Value = 3442824
$ws.Range("C53").Value = 743
$ws.Range("D53").Value = 1109805
$ws.Range("C54").Value = 176
$ws.Range("D54").Value = 260833
$ws.Range("C55").Value = 18
$ws.Range("D55").Value = 27000
$ws.Range("C56").Value = 6613
$ws.Range("D56").Value = 9111676
$ws.Range("C57").Value = 894
$ws.Range("D57").Value = 1312254
$ws.Range("C58").Value = 2243
$ws.Range("D58").Value = 3328113
$ws.Range("C59").Value = 899
$ws.Range("D59").Value = 1338001
$ws.Range("C60").Value = 308
$ws.Range("D60").Value = 461758
$ws.Range("C61").Value = 98
$ws.Range("D61").Value = 147000
$ws.Range("C63").Value = 1309
$ws.Range("D63").Value = 1847385
$ws.Range("C64").Value = 15198
$ws.Range("D64").Value = 21956569
$ws.Range("C65").Value = 44306
$ws.Range("D65").Value = 64848504
$ws.Range("C66").Value = 15590
$ws.Range("D66").Value = 23172214
$ws.Range("C67").Value = 4531
$ws.Range("D67").Value = 6748792
$ws.Range("C68").Value = 907
$ws.Range("D68").Value = 1350096
$ws.Range("C69").Value = 76
$ws.Range("D69").Value = 111330
$ws.Range("C72").Value = 14947
$ws.Range("D72").Value = 19722901
$ws.Range("C73").Value = 50445
$ws.Range("D73").Value = 73419012
$ws.Range("C74").Value = 143772
$ws.Range("D74").Value = 211838266
$ws.Range("C75").Value = 62779
$ws.Range("D75").Value = 93552961
$ws.Range("C76").Value = 20030
$ws.Range("D76").Value = 29928246
$ws.Range("C77").Value = 4714
$ws.Range("D77").Value = 7042723
$ws.Range("C78").Value = 255
$ws.Range("D78").Value = 377670
$ws.Range("C83").Value = 4
$ws.Range("D83").Value = 6000
$ws.Range("C84").Value = 49929
$ws.Range("D84").Value = 68015502
$ws.Range("C85").Value = 4530
$ws.Range("D85").Value = 6563059
$ws.Range("C86").Value = 11397
$ws.Range("D86").Value = 16745098
$ws.Range("C87").Value = 3831
$ws.Range("D87").Value = 5709665
$ws.Range("C88").Value = 1327
$ws.Range("D88").Value = 1982989
$ws.Range("C92").Value = 5290
$ws.Range("D92").Value = 7118424
$ws.Range("C93").Value = 1555
$ws.Range("D93").Value = 2240604
$ws.Range("C94").Value = 5048
$ws.Range("D94").Value = 7436529
$ws.Range("C95").Value = 1913
$ws.Range("D95").Value = 2849946
$ws.Range("C96").Value = 675
$ws.Range("D96").Value = 1011460
$ws.Range("C97").Value = 175
$ws.Range("D97").Value = 261613
$ws.Range("C100").Value = 3449
$ws.Range("D100").Value = 4570529
$ws.Range("C101").Value = 584
$ws.Range("D101").Value = 869664
$ws.Range("C103").Value = 124
$ws.Range("D103").Value = 186000
$ws.Range("C106").Value = 10655
$ws.Range("D106").Value = 15465297
$ws.Range("C107").Value = 28992
$ws.Range("D107").Value = 42603615
$ws.Range("C108").Value = 9710
$ws.Range("D108").Value = 14439828
$ws.Range("C109").Value = 2666
$ws.Range("D109").Value = 3975207
$ws.Range("C110").Value = 484
$ws.Range("D110").Value = 721046
$ws.Range("C113").Value = 9700
$ws.Range("D113").Value = 12824785
$ws.Range("C114").Value = 30078
$ws.Range("D114").Value = 43384967
$ws.Range("C115").Value = 65565
$ws.Range("D115").Value = 95974626
$ws.Range("C116").Value = 21194
$ws.Range("D116").Value = 31497448
$ws.Range("C117").Value = 5990
$ws.Range("D117").Value = 8924826
$ws.Range("C118").Value = 1106
$ws.Range("D118").Value = 1652771
$ws.Range("C119").Value = 76
$ws.Range("D119").Value = 111420
$ws.Range("C122").Value = 4
$ws.Range("D122").Value = 6000
$ws.Range("C123").Value = 25534
$ws.Range("D123").Value = 34131821
$ws.Range("C124").Value = 35484
$ws.Range("D124").Value = 51224956
$ws.Range("C125").Value = 75928
$ws.Range("D125").Value = 111062789
$ws.Range("C126").Value = 23620
$ws.Range("D126").Value = 35061709
$ws.Range("C127").Value = 6329
$ws.Range("D127").Value = 9406004
$ws.Range("C128").Value = 1205
$ws.Range("D128").Value = 1791911
$ws.Range("C132").Value = 31315
$ws.Range("D132").Value = 41611273
$ws.Range("C133").Value = 13108
$ws.Range("D133").Value = 18977379
$ws.Range("C134").Value = 32069
$ws.Range("D134").Value = 47109607
$ws.Range("C135").Value = 11398
$ws.Range("D135").Value = 16936454
$ws.Range("C136").Value = 2930
$ws.Range("D136").Value = 4368504
$ws.Range("C137").Value = 480
$ws.Range("D137").Value = 713990
$ws.Range("C140").Value = 10724
$ws.Range("D140").Value = 14312818
$ws.Range("C141").Value = 34629
$ws.Range("D141").Value = 50018797
$ws.Range("C142").Value = 80399
$ws.Range("D142").Value = 117802265
$ws.Range("C143").Value = 24126
$ws.Range("D143").Value = 35850678
$ws.Range("C144").Value = 6324
$ws.Range("D144").Value = 9437384
$ws.Range("C145").Value = 1417
$ws.Range("D145").Value = 2107730
$ws.Range("C146").Value = 78
$ws.Range("D146").Value = 116630
$ws.Range("C148").Value = 28835
$ws.Range("D148").Value = 38938377
